$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2094","7893843","`"☆ MEDUSA ᴰ☆ᴬ`"","Chinese","4505"),
    @("15883","19940010","maik89","Chinese","3952"),
    @("52180","50627978","`"☆ PANDORA ᴰ☆ᴬ`"","Chinese","2650"),
    @("54312","4154202","`"Klaus Taler`"","Chinese","2608"),
    @("5885","22514681","`"☆ BALOO ᴰ☆ᴬ`"","Chinese","4284"),
    @("14809","33634191","`"Deus Noctem`"","Chinese","3983"),
    @("14833","5653920","SvenZz","Chinese","3982"),
    @("18934","57732675","Grimm","Chinese","3841"),
    @("22892","24195835","`"DSA ToxiX`"","Chinese","3715"),
    @("25462","35493634","`"☆ Thomas S. ᴰ☆ᴬ`"","Chinese","3632"),
    @("25536","5601492","Einherjer1975","Chinese","3629"),
    @("27569","4186901","`"DSA Alex`"","Chinese","3561"),
    @("31646","21042575","Olminator","Chinese","3421"),
    @("34269","58029196","hateboxone","Chinese","3323"),
    @("34316","29621011","SilentStorm","Chinese","3322"),
    @("35271","48588331","Chris","Chinese","3285"),
    @("39831","9121055","JuegoDMierd","Chinese","3074"),
    @("40882","45207062","`"| 123 |`"","Chinese","3025"),
    @("41040","19985443","`"☆ Coura ᴰ☆ᴬ`"","Chinese","3018"),
    @("45589","7315419","Campino1410","Chinese","2840"),
    @("46119","6575464","`"☆ Kinnay ᴰ☆ᴬ`"","Chinese","2822"),
    @("47778","44703584","`"DSA FiftyPark`"","Chinese","2770"),
    @("51416","51102648","Pilsbirne","Chinese","2668"),
    @("55570","51138331","GoldenAquarius","Chinese","2586"),
    @("999999","57949338","`"☆ Neptun01 ᴰ☆ᴬ`"","Chinese","0"),
    @("19676","41468909","To","Chinese","3818"),
    @("23696","7248556","`"☆ Unbroken ᴰ☆ᴬ`"","Chinese","3688"),
    @("24393","58161267","`"☆ Madruk1985 D☆A`"","Chinese","3665"),
    @("26493","44334179","Commando","Chinese","3599"),
    @("27319","14679266","Amboss","Chinese","3570"),
    @("28494","14999725","`"☆ Gandolf ᴰ☆ᴬ`"","Chinese","3531"),
    @("28768","21533545","`"☆ Andy ᴰ☆ᴬ`"","Chinese","3521"),
    @("29674","203024","`"Manuel Dominick`"","Chinese","3491"),
    @("30541","24404822","`"☆ Al Bundy ᴰ☆ᴬ`"","Chinese","3461"),
    @("35357","6010385","`"Jose Javier Davila V`"","Chinese","3281"),
    @("37686","5488187","Vitala93","Chinese","3176"),
    @("39013","3628285","`"Schwatten92 `"","Chinese","3116"),
    @("48144","36985876","`"Eleni Sachinidou`"","Chinese","2759"),
    @("49186","7520688","DSA_AustrianFish","Chinese","2726"),
    @("49226","58992516","Luke","Chinese","2725"),
    @("51926","58187893","MagiendJr","Chinese","2656"),
    @("54103","17695863","Player-17695863","Chinese","2612"),
    @("57722","4539315","Recrush","Chinese","2556"),
    @("60189","58074881","USER39","Chinese","2528"),
    @("62842","6132629","supergamer","Chinese","2504"),
    @("64944","47243338","benja","Chinese","2494"),
    @("68120","56737748","Dennny","Chinese","2331"),
    @("999999","33307536","Player-DsaBruno","Chinese","0"),
    @("999999","17176957","Brutus","Chinese","0")
)

$startRow = 103
$i = 0
foreach ($row in $data) {
    $r = $startRow + $i
    for ($c = 0; $c -lt 5; $c++) {
        $val = $row[$c]
        $cell = $ws.Cells.Item($r, $c + 1)
        if ($val -match '^-?[0-9]+$') {
            $cell.Value = "'" + $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
    $i++
}
